$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell value changes as described by the diff (rows 2-51).
# D-column numeric-looking text values need NumberFormat "@" (Text) first
# so Excel keeps them as text instead of auto-converting to numbers,
# matching the original inlineStr text cells.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.373.98'
$ws.Range("E2").Value = '  -2.64%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.945.35'
$ws.Range("E3").Value = '  -2.58%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '521.04'
$ws.Range("E5").Value = '  -1.32%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.73'
$ws.Range("E6").Value = '  -0.53%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '2.939.15'
$ws.Range("E8").Value = '  -2.41%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.476'
$ws.Range("E9").Value = '  -1.87%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.13'
$ws.Range("E10").Value = '  +1.97%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.145'
$ws.Range("E11").Value = '  -2.07%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.432'
$ws.Range("E12").Value = '  -2.61%  '
$ws.Range("E13").Value = '  -1.99%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.64'
$ws.Range("E14").Value = '  -1.78%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.432.94'
$ws.Range("E15").Value = '  -1.16%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.109'
$ws.Range("E16").Value = '  +0.01%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '60.398.12'
$ws.Range("E17").Value = '  -2.61%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.952.92'
$ws.Range("E18").Value = '  -2.78%  '
$ws.Range("E19").Value = '  -0.31%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '452.00'
$ws.Range("E20").Value = '  -3.69%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.90'
$ws.Range("E21").Value = '  -0.48%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.661'
$ws.Range("E22").Value = '  -3.04%  '
$ws.Range("E23").Value = '  -3.29%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '77.35'
$ws.Range("E24").Value = '  -0.28%  '
$ws.Range("E25").Value = '  -1.31%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.59'
$ws.Range("E27").Value = '  -1.16%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.56'
$ws.Range("E28").Value = '  -5.50%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.998'
$ws.Range("E29").Value = '  -0.17%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '24.98'
$ws.Range("E30").Value = '  -0.85%  '
$ws.Range("E31").Value = '  +2.34%  '
$ws.Range("E32").Value = '  -0.22%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.24'
$ws.Range("E33").Value = '  -4.23%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '54.31'
$ws.Range("E34").Value = '  -2.93%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.23'
$ws.Range("E35").Value = '  +2.72%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.65'
$ws.Range("E36").Value = '  -1.59%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '443.65'
$ws.Range("E37").Value = '  -3.27%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.132.34'
$ws.Range("E38").Value = '  +2.84%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0763'
$ws.Range("E39").Value = '  -0.42%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0375'
$ws.Range("E40").Value = '  -2.03%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.116'
$ws.Range("E41").Value = '  +4.44%  '
$ws.Range("E42").Value = '  +0.05%  '
$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.38'
$ws.Range("E43").Value = '  -4.46%  '
$ws.Range("B44").Value = 'USDe'
$ws.Range("C44").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("E44").Value = '  +0.11%  '
$ws.Range("E45").Value = '  -0.87%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '24.71'
$ws.Range("E46").Value = '  +4.25%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '118.38'
$ws.Range("E47").Value = '  +3.23%  '
$ws.Range("E48").Value = '  +0.59%  '
$ws.Range("E49").Value = '  -3.90%  '
$ws.Range("E50").Value = '  -2.71%  '
$ws.Range("E51").Value = '  +5.95%  '
